$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "System, backup@backdoor.com, system" = "system, backup@backdoor.com, System"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2
    if ($null -ne $val -and $replacements.ContainsKey($val)) {
        $cell.Value = $replacements[$val]
    }
}
